$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skosmos instances")

# Remove whitespace from the "HTW Chur" name -> "HTW-Chur"
$ws.Range("A14").Value = "HTW-Chur"

# Update the active selection shown in the sheet view
$ws.Range("B22").Select()
